# Applies the CORE_holdings update:
#  1) Updates the "as of" date in the confidential disclaimer text (A11)
#     from 2021-05-21 to 2021-05-24.
#  2) Updates the Weight (D2:D8) and Percent Change (E2:E8) values for the
#     fund holdings rows to the new, refreshed figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet is protected - unprotect before editing, then re-apply protection
# with the same options afterwards.
$ws.Unprotect("D382")

# --- 1) Update confidential disclaimer text cell (A11) ---
$ws.Range("A11").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-24 for illustrative purposes only and are subject to change."

# --- 2) Update Weight (D) and Percent Change (E) values ---
$ws.Range("D2").Value = 0.5013101851027676
$ws.Range("E2").Value = 0.005107526881720448

$ws.Range("D3").Value = 0.2439216465473
$ws.Range("E3").Value = 0.01460607848922968

$ws.Range("D4").Value = 0.0951378723604637
$ws.Range("E4").Value = 0.005323868677906018

$ws.Range("D5").Value = 0.1027838251049804
$ws.Range("E5").Value = 0.004749045534965912

$ws.Range("D6").Value = 0.03003117683507748
$ws.Range("E6").Value = 0.00296395448895681

$ws.Range("D7").Value = 0.02681529404941089
$ws.Range("E7").Value = 0.006595995288574796

$ws.Range("D8").Value = 1
$ws.Range("E8").Value = 0.007383705160209253

# --- 3) Restore sheet protection (same options as the original file:
#        contents locked, objects locked, scenarios locked, but column/row
#        formatting allowed) ---
$ws.Protect("D382", $true, $true, $true, $false, $false, $false, $false, $false, $false, $false, $false, $false, $false, $false, $false)
